$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 2
$ws.Range("B10").Value = 2
$ws.Range("B11").Value = 2
$ws.Range("B12").Value = 2

$ws.Range("B13").Select()
